$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("D2", "43.016.15"),
    @("E2", "  -0.66%  "),
    @("D3", "2.341.22"),
    @("E3", "  +0.86%  "),
    @("D4", "1.00"),
    @("E4", "  +0.02%  "),
    @("D5", "303.03"),
    @("E5", "  +0.25%  "),
    @("D6", "94.69"),
    @("E6", "  -3.28%  "),
    @("E7", "  -0.93%  "),
    @("E8", "  +0.04%  "),
    @("E9", "  -1.29%  "),
    @("D10", "34.10"),
    @("E10", "  -4.05%  "),
    @("D11", "0.0784"),
    @("E11", "  -1.51%  "),
    @("D12", "18.65"),
    @("E12", "  -4.63%  "),
    @("E13", "  +1.47%  "),
    @("E14", "  -2.22%  "),
    @("D15", "2.704.93"),
    @("E15", "  +0.59%  "),
    @("D16", "2.352.54"),
    @("E16", "  +1.01%  "),
    @("E17", "  +0.60%  "),
    @("D18", "42.959.92"),
    @("E18", "  -0.61%  "),
    @("D19", "12.09"),
    @("E19", "  -4.84%  "),
    @("E20", "  +1.97%  "),
    @("D21", "0.0₃0889"),
    @("E21", "  -1.08%  "),
    @("D22", "67.92"),
    @("E22", "  -0.12%  "),
    @("D23", "235.95"),
    @("E23", "  -0.55%  "),
    @("E24", "  -1.70%  "),
    @("E25", "  +0.09%  "),
    @("E26", "  -1.82%  "),
    @("D27", "24.63"),
    @("E27", "  -1.71%  "),
    @("D28", "2.23"),
    @("E28", "  +7.60%  "),
    @("D29", "9.16"),
    @("E29", "  +0.40%  "),
    @("D30", "31.44"),
    @("E30", "  -4.90%  "),
    @("D31", "1.00"),
    @("E31", "  -0.04%  "),
    @("E32", "  -0.13%  "),
    @("D33", "0.0737"),
    @("E33", "  +4.59%  "),
    @("D34", "17.27"),
    @("E34", "  -3.51%  "),
    @("E35", "  -2.94%  "),
    @("D36", "1.82"),
    @("E36", "  +2.44%  "),
    @("D37", "2.32"),
    @("E37", "  -1.44%  "),
    @("B38", "Kaspa"),
    @("C38", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"),
    @("D38", "0.101"),
    @("E38", "  -0.50%  "),
    @("B39", "Monero"),
    @("C39", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"),
    @("D39", "122.21"),
    @("E39", "  -25.54%  "),
    @("D40", "2.76"),
    @("E40", "  -1.01%  "),
    @("D41", "22.31"),
    @("E41", "  +17.01%  "),
    @("E42", "  -1.22%  "),
    @("D43", "1.936.40"),
    @("E44", "  +0.18%  "),
    @("D45", "10.11"),
    @("E45", "  -4.89%  "),
    @("D46", "2.10"),
    @("E46", "  +1.24%  "),
    @("D47", "2.72"),
    @("E47", "  -2.80%  "),
    @("D48", "2.571.24"),
    @("E48", "  +0.63%  "),
    @("E49", "  +0.43%  "),
    @("D50", "52.86"),
    @("E50", "  -1.91%  "),
    @("D51", "71.69"),
    @("E51", "  -1.30%  ")
)

foreach ($item in $data) {
    $cellRef = $item[0]
    $val = $item[1]
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $val
}
